$wb = $excel.ActiveWorkbook

$wsShort = $wb.Worksheets.Item("short term")

# --- New TODO items in the "short term" sheet (rows 16-18) ---

# Row 16: new item "70. ..." in A16 (new highlighted style), keep "A" in B16
$wsShort.Range("A16").Value = "70. look at the problem of format of the last page of the pdf file with the argument BMDoutput (cf. ex. transcripto 2018)"
$wsShort.Range("A16").Interior.Color = 49407
$wsShort.Range("A16").Borders.LineStyle = 1
$wsShort.Range("A16").WrapText = $true
$wsShort.Range("B16").Value = "A"

# Row 17: new item "71. ..." in A17 (same new highlighted style), B17 left blank
$wsShort.Range("A17").Value = "71. put an example to help the formating of data especially in a R object"
$wsShort.Range("A17").Interior.Color = 49407
$wsShort.Range("A17").Borders.LineStyle = 1
$wsShort.Range("A17").WrapText = $true

# Row 18: new item "72. ..." in A18 (existing yellow highlighted style), "ML" in B18
$wsShort.Range("A18").Value = "72. Add a function associated with sensitivityplot to get the numrical summaries and to do other plots (boxplots)"
$wsShort.Range("A18").Interior.Color = 65535
$wsShort.Range("A18").Borders.LineStyle = 1
$wsShort.Range("A18").WrapText = $true
$wsShort.Range("B18").Value = "ML"

# --- Active sheet / selection changes ---
# Move the active tab from "done" to "short term", and update the
# selection on "short term" to A21 (leaving "done" selection untouched).
$wsShort.Activate()
$wsShort.Range("A21").Select() | Out-Null
